$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 35133
$ws.Range("B2").Value = "Juliana Barros"
$ws.Range("C2").Value = "P&D"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45092
$ws.Range("G2").Value = 4797.21

# Row 3
$ws.Range("A3").Value = 86894
$ws.Range("B3").Value = "Dr. Davi Lucca Sales"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 4431.95

# Row 4
$ws.Range("A4").Value = 28258
$ws.Range("B4").Value = "Sr. Nicolas Cardoso"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45096
$ws.Range("G4").Value = 7584.72

# Row 5
$ws.Range("A5").Value = 40306
$ws.Range("B5").Value = "Dra. Isis Pereira"
$ws.Range("C5").Value = "Operações"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45095
$ws.Range("G5").Value = 10658.41

# Row 6
$ws.Range("A6").Value = 18893
$ws.Range("B6").Value = "Emilly Freitas"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45081
$ws.Range("G6").Value = 7431.25

# Row 7
$ws.Range("A7").Value = 71892
$ws.Range("B7").Value = "Dra. Alice Silveira"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45078
$ws.Range("G7").Value = 11077.67

# Row 8
$ws.Range("A8").Value = 84452
$ws.Range("B8").Value = "Pietra Cavalcanti"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45081
$ws.Range("G8").Value = 6514.95

# Row 9
$ws.Range("A9").Value = 74908
$ws.Range("B9").Value = "Isaac Viana"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45087
$ws.Range("G9").Value = 7813.41

# Row 10
$ws.Range("A10").Value = 23947
$ws.Range("B10").Value = "João Guilherme Rocha"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 9674.27

# Row 11
$ws.Range("A11").Value = 97929
$ws.Range("B11").Value = "Augusto Alves"
$ws.Range("C11").Value = "Engenharia"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45100
$ws.Range("G11").Value = 4264.97
